$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.002.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.903.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7514'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3079'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.53'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06903'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08019'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7571'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.890.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.266'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.183'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.009.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007761'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.68%  '
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.149.63'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.074'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.322'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1280'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.064'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.346'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.528'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.311'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.053'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05375'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.287'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7384'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.718'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01947'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.764'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.240'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4463'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.949'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8313'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.859'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.057.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05984'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
